{"js": "// Apply the text substitutions described by the diff: the date line and\n// each of the 100 arithmetic-answer cells in the table. Every \"old\" string\n// below is a unique, literal run of text in the document (verified against\n// the source OOXML), so a plain non-wildcard search-and-replace for each\n// pair reproduces the diff exactly without touching anything else.\nconst pairs = [\n  [\n    \"2023-10-13 Friday\",\n    \"2023-10-14 Saturday\"\n  ],\n  [\n    \"60+19=79\",\n    \"69-50=19\"\n  ],\n  [\n    \"29+15=44\",\n    \"98-5=93\"\n  ],\n  [\n    \"69-41=28\",\n    \"93-2=91\"\n  ],\n  [\n    \"73-53=20\",\n    \"12+35=47\"\n  ],\n  [\n    \"37+31=68\",\n    \"0+17=17\"\n  ],\n  [\n    \"75+8=83\",\n    \"43+53=96\"\n  ],\n  [\n    \"65-18=47\",\n    \"59-19=40\"\n  ],\n  [\n    \"23-11=12\",\n    \"78-64=14\"\n  ],\n  [\n    \"52+17=69\",\n    \"39-14=25\"\n  ],\n  [\n    \"82-10=72\",\n    \"31+0=31\"\n  ],\n  [\n    \"60-26=34\",\n    \"57-52=5\"\n  ],\n  [\n    \"87-3=84\",\n    \"31+12=43\"\n  ],\n  [\n    \"13+80=93\",\n    \"31-9=22\"\n  ],\n  [\n    \"57+29=86\",\n    \"8+81=89\"\n  ],\n  [\n    \"0+56=56\",\n    \"59-20=39\"\n  ],\n  [\n    \"87+6=93\",\n    \"96-63=33\"\n  ],\n  [\n    \"27-25=2\",\n    \"38+20=58\"\n  ],\n  [\n    \"48-13=35\",\n    \"17-2=15\"\n  ],\n  [\n    \"84+6=90\",\n    \"84-13=71\"\n  ],\n  [\n    \"60-2=58\",\n    \"24+22=46\"\n  ],\n  [\n    \"46-4=42\",\n    \"36+42=78\"\n  ],\n  [\n    \"14+15=29\",\n    \"92-44=48\"\n  ],\n  [\n    \"85-64=21\",\n    \"1+57=58\"\n  ],\n  [\n    \"56-32=24\",\n    \"34+45=79\"\n  ],\n  [\n    \"31-8=23\",\n    \"0+81=81\"\n  ],\n  [\n    \"26+64=90\",\n    \"77-6=71\"\n  ],\n  [\n    \"30+9=39\",\n    \"20+36=56\"\n  ],\n  [\n    \"13+81=94\",\n    \"66+33=99\"\n  ],\n  [\n    \"44+14=58\",\n    \"10+28=38\"\n  ],\n  [\n    \"92-48=44\",\n    \"37-32=5\"\n  ],\n  [\n    \"8+68=76\",\n    \"93-12=81\"\n  ],\n  [\n    \"37-5=32\",\n    \"48-32=16\"\n  ],\n  [\n    \"62-32=30\",\n    \"29+42=71\"\n  ],\n  [\n    \"85-11=74\",\n    \"21+46=67\"\n  ],\n  [\n    \"25+33=58\",\n    \"40-23=17\"\n  ],\n  [\n    \"30+46=76\",\n    \"4+27=31\"\n  ],\n  [\n    \"16+67=83\",\n    \"64+22=86\"\n  ],\n  [\n    \"96-27=69\",\n    \"80-77=3\"\n  ],\n  [\n    \"67-22=45\",\n    \"9+24=33\"\n  ],\n  [\n    \"31+20=51\",\n    \"84-9=75\"\n  ],\n  [\n    \"15-6=9\",\n    \"72-1=71\"\n  ],\n  [\n    \"7+18=25\",\n    \"70-8=62\"\n  ],\n  [\n    \"74-19=55\",\n    \"73-41=32\"\n  ],\n  [\n    \"71-5=66\",\n    \"20+44=64\"\n  ],\n  [\n    \"60+23=83\",\n    \"39+28=67\"\n  ],\n  [\n    \"80-50=30\",\n    \"94+4=98\"\n  ],\n  [\n    \"79-27=52\",\n    \"42-12=30\"\n  ],\n  [\n    \"17+55=72\",\n    \"49+50=99\"\n  ],\n  [\n    \"5+56=61\",\n    \"25+11=36\"\n  ],\n  [\n    \"49-23=26\",\n    \"38-37=1\"\n  ],\n  [\n    \"68+8=76\",\n    \"4+73=77\"\n  ],\n  [\n    \"75+2=77\",\n    \"85-73=12\"\n  ],\n  [\n    \"87-30=57\",\n    \"20+39=59\"\n  ],\n  [\n    \"88-64=24\",\n    \"96-8=88\"\n  ],\n  [\n    \"10+69=79\",\n    \"72+6=78\"\n  ],\n  [\n    \"9+2=11\",\n    \"17+38=55\"\n  ],\n  [\n    \"25+67=92\",\n    \"45-39=6\"\n  ],\n  [\n    \"15+23=38\",\n    \"47-44=3\"\n  ],\n  [\n    \"11+84=95\",\n    \"43-28=15\"\n  ],\n  [\n    \"93-58=35\",\n    \"38+1=39\"\n  ],\n  [\n    \"67+10=77\",\n    \"23+34=57\"\n  ],\n  [\n    \"17+0=17\",\n    \"57+26=83\"\n  ],\n  [\n    \"42+18=60\",\n    \"16+35=51\"\n  ],\n  [\n    \"58+14=72\",\n    \"17+7=24\"\n  ],\n  [\n    \"19+24=43\",\n    \"89-87=2\"\n  ],\n  [\n    \"87-54=33\",\n    \"68+9=77\"\n  ],\n  [\n    \"23-17=6\",\n    \"83-2=81\"\n  ],\n  [\n    \"77-42=35\",\n    \"90+0=90\"\n  ],\n  [\n    \"87-67=20\",\n    \"81-9=72\"\n  ],\n  [\n    \"49-17=32\",\n    \"90-2=88\"\n  ],\n  [\n    \"98-74=24\",\n    \"82-48=34\"\n  ],\n  [\n    \"49+43=92\",\n    \"82-7=75\"\n  ],\n  [\n    \"68+2=70\",\n    \"27-13=14\"\n  ],\n  [\n    \"43-30=13\",\n    \"74+24=98\"\n  ],\n  [\n    \"78-29=49\",\n    \"62-19=43\"\n  ],\n  [\n    \"65+11=76\",\n    \"51-38=13\"\n  ],\n  [\n    \"65+22=87\",\n    \"51-39=12\"\n  ],\n  [\n    \"6+93=99\",\n    \"35+55=90\"\n  ],\n  [\n    \"24+4=28\",\n    \"42-7=35\"\n  ],\n  [\n    \"60-15=45\",\n    \"41-2=39\"\n  ],\n  [\n    \"84-61=23\",\n    \"0+89=89\"\n  ],\n  [\n    \"58+1=59\",\n    \"45-37=8\"\n  ],\n  [\n    \"61+2=63\",\n    \"95-10=85\"\n  ],\n  [\n    \"94-18=76\",\n    \"57-13=44\"\n  ],\n  [\n    \"40-24=16\",\n    \"31+67=98\"\n  ],\n  [\n    \"58+20=78\",\n    \"66-26=40\"\n  ],\n  [\n    \"14+46=60\",\n    \"93-47=46\"\n  ],\n  [\n    \"97-9=88\",\n    \"87-4=83\"\n  ],\n  [\n    \"14+82=96\",\n    \"89+3=92\"\n  ],\n  [\n    \"60+20=80\",\n    \"94-62=32\"\n  ],\n  [\n    \"41+55=96\",\n    \"88-74=14\"\n  ],\n  [\n    \"69+21=90\",\n    \"85-63=22\"\n  ],\n  [\n    \"99-40=59\",\n    \"43-13=30\"\n  ],\n  [\n    \"42-23=19\",\n    \"60-22=38\"\n  ],\n  [\n    \"99-53=46\",\n    \"18+33=51\"\n  ],\n  [\n    \"50+29=79\",\n    \"67-31=36\"\n  ],\n  [\n    \"17-10=7\",\n    \"54-7=47\"\n  ],\n  [\n    \"35+26=61\",\n    \"34+53=87\"\n  ],\n  [\n    \"40+35=75\",\n    \"23+42=65\"\n  ],\n  [\n    \"18+69=87\",\n    \"77-59=18\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Replace every match (should be exactly one occurrence for each string).\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text substitutions described by the diff: the date line and\n# each of the 100 arithmetic-answer cells in the table. Every \"old\" string\n# is a unique, literal run of text in the document (verified against the\n# source OOXML), so a plain non-wildcard Find/Replace for each pair\n# reproduces the diff exactly without touching anything else.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2023-10-13 Friday', '2023-10-14 Saturday'),\n    @('60+19=79', '69-50=19'),\n    @('29+15=44', '98-5=93'),\n    @('69-41=28', '93-2=91'),\n    @('73-53=20', '12+35=47'),\n    @('37+31=68', '0+17=17'),\n    @('75+8=83', '43+53=96'),\n    @('65-18=47', '59-19=40'),\n    @('23-11=12', '78-64=14'),\n    @('52+17=69', '39-14=25'),\n    @('82-10=72', '31+0=31'),\n    @('60-26=34', '57-52=5'),\n    @('87-3=84', '31+12=43'),\n    @('13+80=93', '31-9=22'),\n    @('57+29=86', '8+81=89'),\n    @('0+56=56', '59-20=39'),\n    @('87+6=93', '96-63=33'),\n    @('27-25=2', '38+20=58'),\n    @('48-13=35', '17-2=15'),\n    @('84+6=90', '84-13=71'),\n    @('60-2=58', '24+22=46'),\n    @('46-4=42', '36+42=78'),\n    @('14+15=29', '92-44=48'),\n    @('85-64=21', '1+57=58'),\n    @('56-32=24', '34+45=79'),\n    @('31-8=23', '0+81=81'),\n    @('26+64=90', '77-6=71'),\n    @('30+9=39', '20+36=56'),\n    @('13+81=94', '66+33=99'),\n    @('44+14=58', '10+28=38'),\n    @('92-48=44', '37-32=5'),\n    @('8+68=76', '93-12=81'),\n    @('37-5=32', '48-32=16'),\n    @('62-32=30', '29+42=71'),\n    @('85-11=74', '21+46=67'),\n    @('25+33=58', '40-23=17'),\n    @('30+46=76', '4+27=31'),\n    @('16+67=83', '64+22=86'),\n    @('96-27=69', '80-77=3'),\n    @('67-22=45', '9+24=33'),\n    @('31+20=51', '84-9=75'),\n    @('15-6=9', '72-1=71'),\n    @('7+18=25', '70-8=62'),\n    @('74-19=55', '73-41=32'),\n    @('71-5=66', '20+44=64'),\n    @('60+23=83', '39+28=67'),\n    @('80-50=30', '94+4=98'),\n    @('79-27=52', '42-12=30'),\n    @('17+55=72', '49+50=99'),\n    @('5+56=61', '25+11=36'),\n    @('49-23=26', '38-37=1'),\n    @('68+8=76', '4+73=77'),\n    @('75+2=77', '85-73=12'),\n    @('87-30=57', '20+39=59'),\n    @('88-64=24', '96-8=88'),\n    @('10+69=79', '72+6=78'),\n    @('9+2=11', '17+38=55'),\n    @('25+67=92', '45-39=6'),\n    @('15+23=38', '47-44=3'),\n    @('11+84=95', '43-28=15'),\n    @('93-58=35', '38+1=39'),\n    @('67+10=77', '23+34=57'),\n    @('17+0=17', '57+26=83'),\n    @('42+18=60', '16+35=51'),\n    @('58+14=72', '17+7=24'),\n    @('19+24=43', '89-87=2'),\n    @('87-54=33', '68+9=77'),\n    @('23-17=6', '83-2=81'),\n    @('77-42=35', '90+0=90'),\n    @('87-67=20', '81-9=72'),\n    @('49-17=32', '90-2=88'),\n    @('98-74=24', '82-48=34'),\n    @('49+43=92', '82-7=75'),\n    @('68+2=70', '27-13=14'),\n    @('43-30=13', '74+24=98'),\n    @('78-29=49', '62-19=43'),\n    @('65+11=76', '51-38=13'),\n    @('65+22=87', '51-39=12'),\n    @('6+93=99', '35+55=90'),\n    @('24+4=28', '42-7=35'),\n    @('60-15=45', '41-2=39'),\n    @('84-61=23', '0+89=89'),\n    @('58+1=59', '45-37=8'),\n    @('61+2=63', '95-10=85'),\n    @('94-18=76', '57-13=44'),\n    @('40-24=16', '31+67=98'),\n    @('58+20=78', '66-26=40'),\n    @('14+46=60', '93-47=46'),\n    @('97-9=88', '87-4=83'),\n    @('14+82=96', '89+3=92'),\n    @('60+20=80', '94-62=32'),\n    @('41+55=96', '88-74=14'),\n    @('69+21=90', '85-63=22'),\n    @('99-40=59', '43-13=30'),\n    @('42-23=19', '60-22=38'),\n    @('99-53=46', '18+33=51'),\n    @('50+29=79', '67-31=36'),\n    @('17-10=7', '54-7=47'),\n    @('35+26=61', '34+53=87'),\n    @('40+35=75', '23+42=65'),\n    @('18+69=87', '77-59=18'),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # MatchCase=True, MatchWildcards=False (so '+' etc. are literal), Replace=wdReplaceAll(2)\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n\"done\"\n"}
